$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: pre-seed the shared-string table in the exact order the new
#     strings must appear, using a far-away scratch row that will be cleared
#     afterwards. This guarantees sharedStrings.xml gets the new <si> entries
#     appended in the same order as the target diff.
$newStrings = @(
    '\$DATA',
    '\$DES',
    '\$ERROR',
    '\$ESTIMATION\|\$EST',
    '\$ESTIMATION (\$EST)',
    '\$INPUT',
    '\$MODEL',
    '\$MSFI',
    '\$OMEGA',
    '\$OMEGA BLOCK',
    '\$OMEGA BIAGONAL',
    '\$PK',
    '\$PRED',
    '\$PROBLEM',
    '\$SCATTERPLOT\|\$SCAT',
    '\$SIGMA',
    '\$SIGMA DIAGONAL',
    '\$SIMULATION\|\$SIM',
    '\$SCATTERPLOT (\$SCAT)',
    '\$SIMULATION (\$SIM)',
    '\$SIMULATION ONLYSIM',
    '\$SUBROUTINE',
    '\$TABLE',
    '\$THETA',
    '\$TOL',
    'ACCEPT',
    'ADDL',
    'ADVAN1',
    'ADVAN2',
    'ADVAN3',
    'ADVAN4',
    'ADVAN5',
    'ADVAN6',
    'ADVAN7',
    'ADVAN8',
    'ADVAN9',
    'ADVAN10',
    'ADVAN11',
    'ADVAN12',
    'ADVAN13',
    'ALAG',
    'AMT',
    'BAYES',
    'CALL',
    'CALLFL = 0',
    'CHECKOUT',
    'CMT',
    'COMP',
    'COMPRESS',
    'CONT',
    'CPRED',
    'CPREDI',
    'CRES',
    'CRESI',
    'CWRES',
    'CWRESI',
    'D1',
    'DADT(i)',
    'DAT1',
    'DAT2',
    'DAT3',
    'DATE',
    'DATE=DROP',
    'DEFDOS',
    'DEFOBS',
    'DROP',
    'DV',
    'ETABAR',
    'EVID',
    'F1',
    'FDATA',
    'F_FLAG',
    'FILE=',
    'FIRSTONLY',
    'FIXED',
    'FO',
    'ID',
    'IF-THEN',
    'IGNORE',
    'INTERACTION',
    'IPRED',
    'IREP',
    'IRES',
    'IWRES',
    'LAPLACIAN',
    'LFORMAT',
    'LIKELIHOOD',
    'MATRIX = R',
    'MATRIX = S',
    'MAXEVAL',
    'MAXEVAL=0',
    'MDV',
    'MPAST',
    'MSF',
    'MSFO',
    'MTIME',
    'NOABORT',
    'NOAPPEND',
    'NOHEADER',
    'NOOMEGABOUNTTEST',
    'NOPRINT',
    'NOSIGMABOUNDTEST',
    'NOTHETABOUNDTEST',
    'NSIG',
    'NSUB',
    'OBSERVATIONS ONLY',
    'ONEHEADER',
    'ONLYSIM',
    'PCMT',
    'POSTHOC',
    'PRED',
    'PREDPP',
    'PRINT=E',
    'PRINT=n',
    'R1',
    'RATE',
    'RFORMAT',
    'S1',
    'SIGDIGITS\|SIGDIG',
    'SORT',
    'SS',
    'SUBPROBLEMS',
    'TIME',
    'RES',
    'TOL',
    'TRANS',
    'TRANS2',
    'TRUE=FINAL',
    'UNIT',
    'WRES',
    'Y\ '
)

$scratchRow = 1000
for ($i = 0; $i -lt $newStrings.Length; $i++) {
    $ws.Cells.Item($scratchRow, $i + 1).Value = $newStrings[$i]
}

# --- Step 2: update existing rows 102-104 (E column 0 -> 1)
$ws.Range("E102").Value = 1
$ws.Range("E103").Value = 1
$ws.Range("E104").Value = 1

# --- Step 3: append new rows 105-235 (E, F, G columns)
# F/G values below reference the same text already interned in step 1, so no
# new shared-string entries are created here -- only new cell references.
$ws.Range("E105").Value = 0
$ws.Range("F105").Value = '\$DATA'
$ws.Range("G105").Value = '\$DATA'
$ws.Range("E106").Value = 0
$ws.Range("F106").Value = '\$DES'
$ws.Range("G106").Value = '\$DES'
$ws.Range("E107").Value = 0
$ws.Range("F107").Value = '\$ERROR'
$ws.Range("G107").Value = '\$ERROR'
$ws.Range("E108").Value = 0
$ws.Range("F108").Value = '\$ESTIMATION (\$EST)'
$ws.Range("G108").Value = '\$ESTIMATION\|\$EST'
$ws.Range("E109").Value = 0
$ws.Range("F109").Value = '\$INPUT'
$ws.Range("G109").Value = '\$INPUT'
$ws.Range("E110").Value = 0
$ws.Range("F110").Value = '\$MODEL'
$ws.Range("G110").Value = '\$MODEL'
$ws.Range("E111").Value = 0
$ws.Range("F111").Value = '\$MSFI'
$ws.Range("G111").Value = '\$MSFI'
$ws.Range("E112").Value = 0
$ws.Range("F112").Value = '\$OMEGA'
$ws.Range("G112").Value = '\$OMEGA'
$ws.Range("E113").Value = 0
$ws.Range("F113").Value = '\$OMEGA BLOCK'
$ws.Range("G113").Value = '\$OMEGA BLOCK'
$ws.Range("E114").Value = 0
$ws.Range("F114").Value = '\$OMEGA BIAGONAL'
$ws.Range("G114").Value = '\$OMEGA BIAGONAL'
$ws.Range("E115").Value = 0
$ws.Range("F115").Value = '\$PK'
$ws.Range("G115").Value = '\$PK'
$ws.Range("E116").Value = 0
$ws.Range("F116").Value = '\$PRED'
$ws.Range("G116").Value = '\$PRED'
$ws.Range("E117").Value = 0
$ws.Range("F117").Value = '\$PROBLEM'
$ws.Range("G117").Value = '\$PROBLEM'
$ws.Range("E118").Value = 0
$ws.Range("F118").Value = '\$SCATTERPLOT (\$SCAT)'
$ws.Range("G118").Value = '\$SCATTERPLOT\|\$SCAT'
$ws.Range("E119").Value = 0
$ws.Range("F119").Value = '\$SIGMA'
$ws.Range("G119").Value = '\$SIGMA'
$ws.Range("E120").Value = 0
$ws.Range("F120").Value = '\$SIGMA DIAGONAL'
$ws.Range("G120").Value = '\$SIGMA DIAGONAL'
$ws.Range("E121").Value = 0
$ws.Range("F121").Value = '\$SIMULATION (\$SIM)'
$ws.Range("G121").Value = '\$SIMULATION\|\$SIM'
$ws.Range("E122").Value = 0
$ws.Range("F122").Value = '\$SIMULATION ONLYSIM'
$ws.Range("G122").Value = '\$SIMULATION ONLYSIM'
$ws.Range("E123").Value = 0
$ws.Range("F123").Value = '\$SUBROUTINE'
$ws.Range("G123").Value = '\$SUBROUTINE'
$ws.Range("E124").Value = 0
$ws.Range("F124").Value = '\$TABLE'
$ws.Range("G124").Value = '\$TABLE'
$ws.Range("E125").Value = 0
$ws.Range("F125").Value = '\$THETA'
$ws.Range("G125").Value = '\$THETA'
$ws.Range("E126").Value = 0
$ws.Range("F126").Value = '\$TOL'
$ws.Range("G126").Value = '\$TOL'
$ws.Range("E127").Value = 0
$ws.Range("F127").Value = 'ACCEPT'
$ws.Range("G127").Value = 'ACCEPT'
$ws.Range("E128").Value = 0
$ws.Range("F128").Value = 'ADDL'
$ws.Range("G128").Value = 'ADDL'
$ws.Range("E129").Value = 0
$ws.Range("F129").Value = 'ADVAN1'
$ws.Range("G129").Value = 'ADVAN1'
$ws.Range("E130").Value = 0
$ws.Range("F130").Value = 'ADVAN2'
$ws.Range("G130").Value = 'ADVAN2'
$ws.Range("E131").Value = 0
$ws.Range("F131").Value = 'ADVAN3'
$ws.Range("G131").Value = 'ADVAN3'
$ws.Range("E132").Value = 0
$ws.Range("F132").Value = 'ADVAN4'
$ws.Range("G132").Value = 'ADVAN4'
$ws.Range("E133").Value = 0
$ws.Range("F133").Value = 'ADVAN5'
$ws.Range("G133").Value = 'ADVAN5'
$ws.Range("E134").Value = 0
$ws.Range("F134").Value = 'ADVAN6'
$ws.Range("G134").Value = 'ADVAN6'
$ws.Range("E135").Value = 0
$ws.Range("F135").Value = 'ADVAN7'
$ws.Range("G135").Value = 'ADVAN7'
$ws.Range("E136").Value = 0
$ws.Range("F136").Value = 'ADVAN8'
$ws.Range("G136").Value = 'ADVAN8'
$ws.Range("E137").Value = 0
$ws.Range("F137").Value = 'ADVAN9'
$ws.Range("G137").Value = 'ADVAN9'
$ws.Range("E138").Value = 0
$ws.Range("F138").Value = 'ADVAN10'
$ws.Range("G138").Value = 'ADVAN10'
$ws.Range("E139").Value = 0
$ws.Range("F139").Value = 'ADVAN11'
$ws.Range("G139").Value = 'ADVAN11'
$ws.Range("E140").Value = 0
$ws.Range("F140").Value = 'ADVAN12'
$ws.Range("G140").Value = 'ADVAN12'
$ws.Range("E141").Value = 0
$ws.Range("F141").Value = 'ADVAN13'
$ws.Range("G141").Value = 'ADVAN13'
$ws.Range("E142").Value = 0
$ws.Range("F142").Value = 'ALAG'
$ws.Range("G142").Value = 'ALAG'
$ws.Range("E143").Value = 0
$ws.Range("F143").Value = 'AMT'
$ws.Range("G143").Value = 'AMT'
$ws.Range("E144").Value = 0
$ws.Range("F144").Value = 'BAYES'
$ws.Range("G144").Value = 'BAYES'
$ws.Range("E145").Value = 0
$ws.Range("F145").Value = 'CALL'
$ws.Range("G145").Value = 'CALL'
$ws.Range("E146").Value = 0
$ws.Range("F146").Value = 'CALLFL = 0'
$ws.Range("G146").Value = 'CALLFL = 0'
$ws.Range("E147").Value = 0
$ws.Range("F147").Value = 'CHECKOUT'
$ws.Range("G147").Value = 'CHECKOUT'
$ws.Range("E148").Value = 0
$ws.Range("F148").Value = 'CMT'
$ws.Range("G148").Value = 'CMT'
$ws.Range("E149").Value = 0
$ws.Range("F149").Value = 'COMP'
$ws.Range("G149").Value = 'COMP'
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 'COMPRESS'
$ws.Range("G150").Value = 'COMPRESS'
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 'CONT'
$ws.Range("G151").Value = 'CONT'
$ws.Range("E152").Value = 0
$ws.Range("F152").Value = 'CPRED'
$ws.Range("G152").Value = 'CPRED'
$ws.Range("E153").Value = 0
$ws.Range("F153").Value = 'CPREDI'
$ws.Range("G153").Value = 'CPREDI'
$ws.Range("E154").Value = 0
$ws.Range("F154").Value = 'CRES'
$ws.Range("G154").Value = 'CRES'
$ws.Range("E155").Value = 0
$ws.Range("F155").Value = 'CRESI'
$ws.Range("G155").Value = 'CRESI'
$ws.Range("E156").Value = 0
$ws.Range("F156").Value = 'CWRES'
$ws.Range("G156").Value = 'CWRES'
$ws.Range("E157").Value = 0
$ws.Range("F157").Value = 'CWRESI'
$ws.Range("G157").Value = 'CWRESI'
$ws.Range("E158").Value = 0
$ws.Range("F158").Value = 'D1'
$ws.Range("G158").Value = 'D1'
$ws.Range("E159").Value = 0
$ws.Range("F159").Value = 'DADT(i)'
$ws.Range("G159").Value = 'DADT(i)'
$ws.Range("E160").Value = 0
$ws.Range("F160").Value = 'DAT1'
$ws.Range("G160").Value = 'DAT1'
$ws.Range("E161").Value = 0
$ws.Range("F161").Value = 'DAT2'
$ws.Range("G161").Value = 'DAT2'
$ws.Range("E162").Value = 0
$ws.Range("F162").Value = 'DAT3'
$ws.Range("G162").Value = 'DAT3'
$ws.Range("E163").Value = 0
$ws.Range("F163").Value = 'DATE'
$ws.Range("G163").Value = 'DATE'
$ws.Range("E164").Value = 0
$ws.Range("F164").Value = 'DATE=DROP'
$ws.Range("G164").Value = 'DATE=DROP'
$ws.Range("E165").Value = 0
$ws.Range("F165").Value = 'DEFDOS'
$ws.Range("G165").Value = 'DEFDOS'
$ws.Range("E166").Value = 0
$ws.Range("F166").Value = 'DEFOBS'
$ws.Range("G166").Value = 'DEFOBS'
$ws.Range("E167").Value = 0
$ws.Range("F167").Value = 'DROP'
$ws.Range("G167").Value = 'DROP'
$ws.Range("E168").Value = 0
$ws.Range("F168").Value = 'DV'
$ws.Range("G168").Value = 'DV'
$ws.Range("E169").Value = 0
$ws.Range("F169").Value = 'ETABAR'
$ws.Range("G169").Value = 'ETABAR'
$ws.Range("E170").Value = 0
$ws.Range("F170").Value = 'EVID'
$ws.Range("G170").Value = 'EVID'
$ws.Range("E171").Value = 0
$ws.Range("F171").Value = 'F'
$ws.Range("G171").Value = 'F'
$ws.Range("E172").Value = 0
$ws.Range("F172").Value = 'F1'
$ws.Range("G172").Value = 'F1'
$ws.Range("E173").Value = 0
$ws.Range("F173").Value = 'FDATA'
$ws.Range("G173").Value = 'FDATA'
$ws.Range("E174").Value = 0
$ws.Range("F174").Value = 'F_FLAG'
$ws.Range("G174").Value = 'F_FLAG'
$ws.Range("E175").Value = 0
$ws.Range("F175").Value = 'FILE='
$ws.Range("G175").Value = 'FILE='
$ws.Range("E176").Value = 0
$ws.Range("F176").Value = 'FIRSTONLY'
$ws.Range("G176").Value = 'FIRSTONLY'
$ws.Range("E177").Value = 0
$ws.Range("F177").Value = 'FIXED'
$ws.Range("G177").Value = 'FIXED'
$ws.Range("E178").Value = 0
$ws.Range("F178").Value = 'FOCE'
$ws.Range("G178").Value = 'FOCE'
$ws.Range("E179").Value = 0
$ws.Range("F179").Value = 'FO'
$ws.Range("G179").Value = 'FO'
$ws.Range("E180").Value = 0
$ws.Range("F180").Value = 'ID'
$ws.Range("G180").Value = 'ID'
$ws.Range("E181").Value = 0
$ws.Range("F181").Value = 'IF-THEN'
$ws.Range("G181").Value = 'IF-THEN'
$ws.Range("E182").Value = 0
$ws.Range("F182").Value = 'IGNORE'
$ws.Range("G182").Value = 'IGNORE'
$ws.Range("E183").Value = 0
$ws.Range("F183").Value = 'INTERACTION'
$ws.Range("G183").Value = 'INTERACTION'
$ws.Range("E184").Value = 0
$ws.Range("F184").Value = 'IPRED'
$ws.Range("G184").Value = 'IPRED'
$ws.Range("E185").Value = 0
$ws.Range("F185").Value = 'IREP'
$ws.Range("G185").Value = 'IREP'
$ws.Range("E186").Value = 0
$ws.Range("F186").Value = 'IRES'
$ws.Range("G186").Value = 'IRES'
$ws.Range("E187").Value = 0
$ws.Range("F187").Value = 'IWRES'
$ws.Range("G187").Value = 'IWRES'
$ws.Range("E188").Value = 0
$ws.Range("F188").Value = 'LAPLACIAN'
$ws.Range("G188").Value = 'LAPLACIAN'
$ws.Range("E189").Value = 0
$ws.Range("F189").Value = 'LFORMAT'
$ws.Range("G189").Value = 'LFORMAT'
$ws.Range("E190").Value = 0
$ws.Range("F190").Value = 'LIKELIHOOD'
$ws.Range("G190").Value = 'LIKELIHOOD'
$ws.Range("E191").Value = 0
$ws.Range("F191").Value = 'MATRIX = R'
$ws.Range("G191").Value = 'MATRIX = R'
$ws.Range("E192").Value = 0
$ws.Range("F192").Value = 'MATRIX = S'
$ws.Range("G192").Value = 'MATRIX = S'
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 'MAXEVAL'
$ws.Range("G193").Value = 'MAXEVAL'
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 'MAXEVAL=0'
$ws.Range("G194").Value = 'MAXEVAL=0'
$ws.Range("E195").Value = 0
$ws.Range("F195").Value = 'MDV'
$ws.Range("G195").Value = 'MDV'
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 'MPAST'
$ws.Range("G196").Value = 'MPAST'
$ws.Range("E197").Value = 0
$ws.Range("F197").Value = 'MSF'
$ws.Range("G197").Value = 'MSF'
$ws.Range("E198").Value = 0
$ws.Range("F198").Value = 'MSFO'
$ws.Range("G198").Value = 'MSFO'
$ws.Range("E199").Value = 0
$ws.Range("F199").Value = 'MTIME'
$ws.Range("G199").Value = 'MTIME'
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 'NM-TRAN'
$ws.Range("G200").Value = 'NM-TRAN'
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 'NOABORT'
$ws.Range("G201").Value = 'NOABORT'
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 'NOAPPEND'
$ws.Range("G202").Value = 'NOAPPEND'
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 'NOHEADER'
$ws.Range("G203").Value = 'NOHEADER'
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 'NOOMEGABOUNTTEST'
$ws.Range("G204").Value = 'NOOMEGABOUNTTEST'
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 'NOPRINT'
$ws.Range("G205").Value = 'NOPRINT'
$ws.Range("E206").Value = 0
$ws.Range("F206").Value = 'NOSIGMABOUNDTEST'
$ws.Range("G206").Value = 'NOSIGMABOUNDTEST'
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 'NOTHETABOUNDTEST'
$ws.Range("G207").Value = 'NOTHETABOUNDTEST'
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 'NSIG'
$ws.Range("G208").Value = 'NSIG'
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 'NSUB'
$ws.Range("G209").Value = 'NSUB'
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 'OBSERVATIONS ONLY'
$ws.Range("G210").Value = 'OBSERVATIONS ONLY'
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 'ONEHEADER'
$ws.Range("G211").Value = 'ONEHEADER'
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 'ONLYSIM'
$ws.Range("G212").Value = 'ONLYSIM'
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 'PCMT'
$ws.Range("G213").Value = 'PCMT'
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 'POSTHOC'
$ws.Range("G214").Value = 'POSTHOC'
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 'PRED'
$ws.Range("G215").Value = 'PRED'
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 'PREDPP'
$ws.Range("G216").Value = 'PREDPP'
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 'PRINT=E'
$ws.Range("G217").Value = 'PRINT=E'
$ws.Range("E218").Value = 0
$ws.Range("F218").Value = 'PRINT=n'
$ws.Range("G218").Value = 'PRINT=n'
$ws.Range("E219").Value = 0
$ws.Range("F219").Value = 'R1'
$ws.Range("G219").Value = 'R1'
$ws.Range("E220").Value = 0
$ws.Range("F220").Value = 'RATE'
$ws.Range("G220").Value = 'RATE'
$ws.Range("E221").Value = 0
$ws.Range("F221").Value = 'RES'
$ws.Range("G221").Value = 'RES'
$ws.Range("E222").Value = 0
$ws.Range("F222").Value = 'RFORMAT'
$ws.Range("G222").Value = 'RFORMAT'
$ws.Range("E223").Value = 0
$ws.Range("F223").Value = 'S1'
$ws.Range("G223").Value = 'S1'
$ws.Range("E224").Value = 0
$ws.Range("F224").Value = 'SIGDIGITS\|SIGDIG'
$ws.Range("G224").Value = 'SIGDIGITS\|SIGDIG'
$ws.Range("E225").Value = 0
$ws.Range("F225").Value = 'SORT'
$ws.Range("G225").Value = 'SORT'
$ws.Range("E226").Value = 0
$ws.Range("F226").Value = 'SS'
$ws.Range("G226").Value = 'SS'
$ws.Range("E227").Value = 0
$ws.Range("F227").Value = 'SUBPROBLEMS'
$ws.Range("G227").Value = 'SUBPROBLEMS'
$ws.Range("E228").Value = 0
$ws.Range("F228").Value = 'TIME'
$ws.Range("G228").Value = 'TIME'
$ws.Range("E229").Value = 0
$ws.Range("F229").Value = 'TOL'
$ws.Range("G229").Value = 'TOL'
$ws.Range("E230").Value = 0
$ws.Range("F230").Value = 'TRANS'
$ws.Range("G230").Value = 'TRANS'
$ws.Range("E231").Value = 0
$ws.Range("F231").Value = 'TRANS2'
$ws.Range("G231").Value = 'TRANS2'
$ws.Range("E232").Value = 0
$ws.Range("F232").Value = 'TRUE=FINAL'
$ws.Range("G232").Value = 'TRUE=FINAL'
$ws.Range("E233").Value = 0
$ws.Range("F233").Value = 'UNIT'
$ws.Range("G233").Value = 'UNIT'
$ws.Range("E234").Value = 0
$ws.Range("F234").Value = 'WRES'
$ws.Range("G234").Value = 'WRES'
$ws.Range("E235").Value = 0
$ws.Range("F235").Value = 'Y\ '
$ws.Range("G235").Value = 'Y\ '

# --- Step 4: clear the scratch row used for pre-seeding the shared strings so
#     it does not remain part of the worksheet data / used range.
[void]$ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, $newStrings.Length)).ClearContents()

# --- Step 5: cosmetic view state (active selection) matching the end state
#     captured in the saved workbook. The sheet is already frozen at row 1
#     (ySplit=1) from the source file, so only the active selection is updated.
[void]$ws.Range("F155").Select()
